$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New forecast values for B1:L19 (weekly forecast window rolled forward one week,
# all model outputs refreshed) per "Beginning of Feb update".
$data = @(
  @(45685, 45692, 45699, 45706, 45713, 45720, 45727, 45734, 45741, 45748, 45755),
  @(42.4, 42.4, 43, 44.1, 45.8, 48, 50.7, 53.8, 57.2, 60.8, 64.40000000000001),
  @(49.3, 50.5, 52.3, 54.5, 57.1, 60.1, 63.4, 67.09999999999999, 71.09999999999999, 75.09999999999999, 78.90000000000001),
  @(34.6, 34.8, 35.5, 36.6, 38.2, 40, 42.3, 44.9, 47.8, 51, 54.3),
  @(34.7, 33.4, 32.7, 32.7, 33.5, 34.9, 37.2, 40.1, 43.6, 47.6, 51.7),
  @(50.9, 50.9, 51.1, 51.5, 52.2, 53.1, 54.5, 56.2, 58.4, 60.9, 63.6),
  @(24, 23.7, 24.4, 25.9, 28.2, 31, 34.1, 37.3, 40.6, 43.7, 46.5),
  @(42.6, 42.4, 42.7, 43.6, 45, 46.9, 49.4, 52.5, 55.9, 59.7, 63.5),
  @(44.8, 45.5, 46.7, 48.5, 50.6, 52.9, 55.5, 58.1, 60.8, 63.4, 65.8),
  @(31.1, 29.9, 29.2, 29.1, 29.5, 30.5, 32, 34.2, 36.8, 39.9, 43.5),
  @(54.5, 54.4, 54.6, 55.2, 56.2, 57.8, 59.9, 62.6, 65.59999999999999, 68.90000000000001, 72.3),
  @(44.6, 44.9, 45.5, 46.5, 47.8, 49.4, 51.3, 53.4, 55.8, 58.3, 60.9),
  @(57.5, 58, 58.5, 59, 59.5, 60.1, 60.9, 62.2, 63.8, 65.90000000000001, 68.3),
  @(45.7, 44.5, 43.8, 43.7, 44.2, 45.4, 47.2, 49.6, 52.3, 55.3, 58.2),
  @(50.5, 50.9, 51.2, 51.4, 51.4, 51.2, 51, 50.8, 50.6, 50.6, 50.7),
  @(34.1, 33.3, 33.4, 34.4, 36.3, 38.9, 42.2, 46, 50.3, 54.9, 59.4),
  @(37.1, 35.4, 34.5, 34.5, 35.5, 37.4, 40.2, 43.7, 47.9, 52.3, 56.8),
  @(45.8, 47.4, 50.5, 55.1, 61, 68, 75.90000000000001, 84.3, 92.90000000000001, 101.3, 109.1),
  @(39.6, 40.9, 43.6, 47.5, 52.5, 58.3, 64.7, 71.5, 78.3, 85, 91.2)
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r + 1, $c + 2).Value = $row[$c]
  }
}

Write-Output "done"
